$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so that numeric-looking
# strings (e.g. "673.42") are not auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '69.651.43'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.705.15'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '673.42'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').Value = '162.01'
$ws.Range('E6').Value = '  +2.48%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').Value = '32.95'
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('D14').Value = '3.694.86'
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').Value = '69.631.86'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('D17').Value = '16.34'
$ws.Range('E17').Value = '  +2.43%  '
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('D19').Value = '474.26'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').Value = '9.81'
$ws.Range('E20').Value = '  -2.26%  '
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('D22').Value = '80.43'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').Value = '3.853.50'
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('D24').Value = '0.0000127'
$ws.Range('E24').Value = '  +5.55%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '11.03'
$ws.Range('E26').Value = '  +1.19%  '
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('E29').Value = '  -0.84%  '
$ws.Range('D30').Value = '2.02'
$ws.Range('E30').Value = '  +1.51%  '
$ws.Range('D31').Value = '0.169'
$ws.Range('E31').Value = '  +7.60%  '
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('D34').Value = '26.96'
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('D35').Value = '3.694.80'
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('E36').Value = '  +4.39%  '
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = '0.0915'
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('D42').Value = '173.92'
$ws.Range('E42').Value = '  +2.98%  '
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('E45').Value = '  +2.22%  '
$ws.Range('D46').Value = '0.000279'
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('E47').Value = '  +2.24%  '
$ws.Range('D48').Value = '27.86'
$ws.Range('E48').Value = '  +3.74%  '
$ws.Range('D49').Value = '1.09'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('E51').Value = '  +0.68%  '

# Restore the original (default) cell formatting now that the text values
# are safely stored, so no stray style/number-format is left behind.
$ws.Range("D2:D51").ClearFormats()
